$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New sample/import data in the header area ---
$ws.Range("B2").Value = "Thai"
$ws.Range("B3").Value = 12

# --- Relabel the "#" / "Grade Average" header cells (new shared strings
#     are appended to the table in this order: Thai, Average, Year, A, M12ff@) ---
$ws.Range("K5").Value = "Average"
$ws.Range("L4").Value = "Year"
$ws.Range("D6").Value = "A"
$ws.Range("B6").Value = "M12ff@"

# --- Relabel the remaining header cells (reuse existing "Grade" string) ---
$ws.Range("K4").Value = "Grade"
$ws.Range("L5").Value = "Grade"

# --- Narrower column K ---
$ws.Columns.Item(11).ColumnWidth = 9.79

# --- New sample data row 6 ---
$ws.Range("C6").Value = 1.5555000000000001
$ws.Range("E6").Value = 1.2
$ws.Range("F6").Value = 4

# --- Hyperlink on B6 (adds the Hyperlink style + relationship) ---
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:M12ff@gmail.com")

# --- Selection moves to B6 ---
[void]$ws.Range("B6").Select()
